$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-5 down to 3-6)
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with data
$ws.Range("A2").Value = "Adidas"
$ws.Range("B2").Value = "m"
$ws.Range("C2").Value = 35.5
$ws.Range("D2").Value = 5.5
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 37
$ws.Range("F2").NumberFormat = "# ?/?"

# Update selection to match target
$ws.Range("C11").Select()
